$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price/volume columns keep text formatting so numeric-looking
# strings (e.g. "1.00") are not coerced into numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "37.367.46"
$ws.Range("E2").Value = "  +2.48%  "

# Row 3
$ws.Range("D3").Value = "2.064.36"
$ws.Range("E3").Value = "  +3.89%  "

# Row 4
$ws.Range("E4").Value = "  +0.06%  "

# Row 5
$ws.Range("D5").Value = "234.66"
$ws.Range("E5").Value = "  -0.21%  "

# Row 6
$ws.Range("D6").Value = "0.614"
$ws.Range("E6").Value = "  +3.10%  "

# Row 7
$ws.Range("D7").Value = "58.34"
$ws.Range("E7").Value = "  +7.17%  "

# Row 8
$ws.Range("E8").Value = "  +0.05%  "

# Row 9
$ws.Range("E9").Value = "  +3.63%  "

# Row 10
$ws.Range("D10").Value = "58.84"
$ws.Range("E10").Value = "  +2.19%  "

# Row 11
$ws.Range("E11").Value = "  +2.00%  "

# Row 12
$ws.Range("E12").Value = "  +3.00%  "

# Row 13
$ws.Range("D13").Value = "2.368.16"
$ws.Range("E13").Value = "  +3.96%  "

# Row 14
$ws.Range("D14").Value = "14.64"
$ws.Range("E14").Value = "  +3.69%  "

# Row 15
$ws.Range("E15").Value = "  +4.21%  "

# Row 16
$ws.Range("D16").Value = "0.777"
$ws.Range("E16").Value = "  +3.06%  "

# Row 17
$ws.Range("E17").Value = "  +2.56%  "

# Row 18
$ws.Range("D18").Value = "2.074.79"
$ws.Range("E18").Value = "  +4.21%  "

# Row 19
$ws.Range("D19").Value = "37.565.61"
$ws.Range("E19").Value = "  +3.26%  "

# Row 20
$ws.Range("D20").Value = "6.14"
$ws.Range("E20").Value = "  +17.28%  "

# Row 21
$ws.Range("E21").Value = "  +1.98%  "

# Row 22
$ws.Range("E22").Value = "  +1.37%  "

# Row 23
$ws.Range("D23").Value = "226.53"

# Row 24
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  -0.07%  "

# Row 25
$ws.Range("D25").Value = "2.43"
$ws.Range("E25").Value = "  +1.66%  "

# Row 26
$ws.Range("E26").Value = "  +1.10%  "

# Row 27
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "164.85"
$ws.Range("E27").Value = "  +1.41%  "

# Row 28
$ws.Range("B28").Value = "ImmutableX"
$ws.Range("C28").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D28").Value = "1.51"
$ws.Range("E28").Value = "  +14.60%  "

# Row 29
$ws.Range("E29").Value = "  +2.14%  "

# Row 30
$ws.Range("D30").Value = "19.22"
$ws.Range("E30").Value = "  +2.19%  "

# Row 31
$ws.Range("E31").Value = "  -1.56%  "

# Row 32
$ws.Range("D32").Value = "0.119"
$ws.Range("E32").Value = "  +2.31%  "

# Row 33
$ws.Range("D33").Value = "4.49"
$ws.Range("E33").Value = "  +3.15%  "

# Row 34
$ws.Range("E34").Value = "  +3.06%  "

# Row 35
$ws.Range("D35").Value = "2.55"
$ws.Range("E35").Value = "  +9.06%  "

# Row 36
$ws.Range("D36").Value = "4.56"
$ws.Range("E36").Value = "  +7.49%  "

# Row 37
$ws.Range("E37").Value = "  +1.47%  "

# Row 38
$ws.Range("E38").Value = "  +0.08%  "

# Row 39
$ws.Range("D39").Value = "1.78"
$ws.Range("E39").Value = "  +0.89%  "

# Row 40
$ws.Range("E40").Value = "  +5.94%  "

# Row 41
$ws.Range("D41").Value = "0.0989"
$ws.Range("E41").Value = "  +5.95%  "

# Row 42
$ws.Range("E42").Value = "  -1.82%  "

# Row 43
$ws.Range("D43").Value = "1.469.15"
$ws.Range("E43").Value = "  +1.07%  "

# Row 44
$ws.Range("D44").Value = "96.21"
$ws.Range("E44").Value = "  +7.97%  "

# Row 45
$ws.Range("D45").Value = "4.34"
$ws.Range("E45").Value = "  +19.89%  "

# Row 46
$ws.Range("E46").Value = "  +6.40%  "

# Row 48
$ws.Range("D48").Value = "15.89"
$ws.Range("E48").Value = "  +5.94%  "

# Row 49
$ws.Range("E49").Value = "  +3.99%  "

# Row 50
$ws.Range("D50").Value = "7.27"
$ws.Range("E50").Value = "  +6.48%  "

# Row 51
$ws.Range("E51").Value = "  +2.22%  "
